$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing city values (validation typos introduced on some cities)
$ws.Range("A2").Value = "London1"
$ws.Range("A3").Value = "Amsterdan"
$ws.Range("A4").Value = "Paris2"
$ws.Range("A5").Value = "Rome"
$ws.Range("A6").Value = "Madrid$$$"

# Add a new blank, underlined cell below the list for input validation
$ws.Range("A7").Font.Underline = $true

# Move selection to the newly added cell
$ws.Range("A7").Select()
